$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15: Técnico Subsequente em Edificações - Campus Congonhas
$ws.Range("E15").Value = 99
$ws.Range("F15").Value = 48
$ws.Range("H15").Value = 59

# Row 16: Técnico Subsequente em Mecânica - Campus Congonhas
$ws.Range("E16").Value = 323
$ws.Range("F16").Value = 91
$ws.Range("H16").Value = 178
